# Adds two new columns, I ("I0") and J ("IF"), to the existing data sheet.
# Header cells I1/J1 get the same formatting as the existing header cells
# (bold, bordered, centered) by copying the format from H1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting already used by the other header cells (B1:H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows (row number, I value, J value) ---
$rowData = @(
    @(2, 8, 8),
    @(3, 8, 8),
    @(4, 8, 8),
    @(5, 8, 8),
    @(6, 8, 8),
    @(7, 7, 8),
    @(8, 7, 7),
    @(9, 8, 8),
    @(10, 7, 8),
    @(11, 8, 8),
    @(12, 8, 8),
    @(13, 7, 8),
    @(14, 8, 8),
    @(15, 8, 8),
    @(16, 7, 7),
    @(17, 8, 8),
    @(18, 6, 7),
    @(19, 8, 8),
    @(20, 7, 8),
    @(21, 8, 8),
    @(22, 8, 8),
    @(23, 9, 9),
    @(24, 9, 9),
    @(25, 8, 8),
    @(26, 7, 7),
    @(27, 7, 7),
    @(28, 8, 8),
    @(29, 7, 7),
    @(30, 8, 8),
    @(31, 6, 6),
    @(32, 8, 8),
    @(33, 8, 8),
    @(34, 7, 7),
    @(35, 7, 7),
    @(36, 7, 7),
    @(37, 7, 7),
    @(38, 7, 7),
    @(39, 9, 9),
    @(40, 7, 7),
    @(41, 10, 10),
    @(42, 7, 7),
    @(43, 6, 6),
    @(44, 6, 6),
    @(45, 9, 9),
    @(46, 8, 8),
    @(47, 6, 6),
    @(48, 7, 7),
    @(49, 8, 8),
    @(50, 6, 7),
    @(51, 7, 8),
    @(52, 6, 7),
    @(53, 8, 8),
    @(54, 7, 8),
    @(55, 8, 8),
    @(56, 5, 6),
    @(57, 7, 7),
    @(58, 6, 6),
    @(59, 4, 5),
    @(60, 7, 7),
    @(61, 7, 7),
    @(62, 11, 11),
    @(63, 8, 8),
    @(64, 9, 9),
    @(65, 7, 7),
    @(66, 7, 7),
    @(67, 8, 8),
    @(68, 6, 6),
    @(69, 7, 7),
    @(70, 8, 8),
    @(71, 6, 6),
    @(72, 6, 6)
)

foreach ($entry in $rowData) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
